$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 221, pushing the existing rows 221:298
# down to 222:299 (carrying their styles/formats with them).
$ws.Rows("221").Insert()

# Populate the newly inserted row 221 with the new data point.
$ws.Range("A221").Value = 11
$ws.Range("B221").Value = "Vega Monumental Concepción"
$ws.Range("C221").Value = "Bíobío"
$ws.Range("D221").Value = 44726
$ws.Range("E221").Value = 8
$ws.Range("F221").Value = 100112023
$ws.Range("G221").Value = "Brócoli"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 2200
$ws.Range("K221").Value = 800
$ws.Range("L221").Value = 900
$ws.Range("M221").Value = 855
$ws.Range("N221").Value = "`$/unidad"
$ws.Range("O221").Value = "Región Metropolitana"
$ws.Range("P221").Value = 855
$ws.Range("Q221").Value = 1
$ws.Range("R221").Value = "Hortaliza"
